$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text (rich-text shared strings): bump the report volume/number and
# shift the covered week forward by 7 days.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/17/2024  Through  6/23/2024"

# ---------------------------------------------------------------------------
# Helper: re-type a cell (number<->text) while keeping it looking like its
# neighbours by pasting the number-format/style from a donor cell that
# already has the desired type.
# ---------------------------------------------------------------------------
function Set-CellAs {
    param($target, $value, $donor)
    $ws.Range($target).Value = $value
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($target).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = 44.444444444444

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
Set-CellAs "C17" 2 "D17"
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 6
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 44
$ws.Range("H17").Value = -50
$ws.Range("K17").Value = 4.545454545454
$ws.Range("L17").Value = 130
$ws.Range("M17").Value = 119.047619047619
$ws.Range("N17").Value = -13.207547169811

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -71.428571428571
$ws.Range("J18").Value = 28
$ws.Range("K18").Value = -46.428571428571
$ws.Range("L18").Value = 7.142857142857
$ws.Range("M18").Value = -73.684210526315
$ws.Range("N18").Value = -89.361702127659

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 33.333333333333
$ws.Range("I19").Value = 141
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = 0.714285714285
$ws.Range("L19").Value = 11.904761904761
$ws.Range("M19").Value = 123.809523809524
$ws.Range("N19").Value = 63.953488372093

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
Set-CellAs "C20" "'0" "D14"
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -57.142857142857
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -57.446808510638
$ws.Range("N20").Value = -94.047619047619

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -58.823529411764
$ws.Range("F21").Value = 42
$ws.Range("G21").Value = 51
$ws.Range("H21").Value = -17.647058823529
$ws.Range("I21").Value = 239
$ws.Range("J21").Value = 263
$ws.Range("K21").Value = -9.125475285171
$ws.Range("L21").Value = 12.206572769953
$ws.Range("M21").Value = 39.766081871345
$ws.Range("N21").Value = -62.830482115085

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 30
$ws.Range("G24").Value = 35
$ws.Range("H24").Value = -14.285714285714
$ws.Range("I24").Value = 194
$ws.Range("J24").Value = 232
$ws.Range("K24").Value = -16.379310344827
$ws.Range("L24").Value = -1.020408163265
$ws.Range("M24").Value = -15.652173913043

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 3
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 15.384615384615
$ws.Range("I25").Value = 95
$ws.Range("J25").Value = 111
$ws.Range("K25").Value = -14.414414414414
$ws.Range("L25").Value = 126.190476190476

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -66.666666666666
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = -27.777777777777
$ws.Range("I26").Value = 77
$ws.Range("J26").Value = 92
$ws.Range("K26").Value = -16.304347826087
$ws.Range("L26").Value = -2.531645569620
$ws.Range("M26").Value = -23

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 1
Set-CellAs "D28" "'0" "M28"
Set-CellAs "E28" "***.*" "N28"
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = -40
$ws.Range("L28").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------------
$ws.Range("L31").Value = -66.666666666666

# ---------------------------------------------------------------------------
# Row 33 - Traffic Fatalities
# ---------------------------------------------------------------------------
Set-CellAs "C33" 1 "J33"
Set-CellAs "F33" 1 "J33"
Set-CellAs "I33" 1 "J33"
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
